# Update "想去人数" (interested-count) figures in the "展览" and "全部类型"
# sheets to reflect the latest scrape (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" --------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 3123
$ws1.Range("F4").Value  = 1077
$ws1.Range("F5").Value  = 76
$ws1.Range("F7").Value  = 268
$ws1.Range("F10").Value = 15470
$ws1.Range("F11").Value = 226
$ws1.Range("F14").Value = 6100
$ws1.Range("F16").Value = 102
$ws1.Range("F21").Value = 26
$ws1.Range("F22").Value = 111
$ws1.Range("F23").Value = 7
$ws1.Range("F27").Value = 853
$ws1.Range("F28").Value = 17
$ws1.Range("F29").Value = 4986
$ws1.Range("F30").Value = 137
$ws1.Range("F31").Value = 10972
$ws1.Range("F36").Value = 3788
$ws1.Range("F37").Value = 259

# --- Sheet "全部类型" ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value  = 3123
$ws4.Range("F5").Value  = 1077
$ws4.Range("F6").Value  = 76
$ws4.Range("F8").Value  = 268
$ws4.Range("F11").Value = 15470
$ws4.Range("F12").Value = 226
$ws4.Range("F15").Value = 6100
$ws4.Range("F17").Value = 102
$ws4.Range("F22").Value = 26
$ws4.Range("F23").Value = 111
$ws4.Range("F24").Value = 7
$ws4.Range("F28").Value = 853
$ws4.Range("F29").Value = 17
$ws4.Range("F30").Value = 4986
$ws4.Range("F31").Value = 137
$ws4.Range("F33").Value = 10972
$ws4.Range("F38").Value = 3788
$ws4.Range("F39").Value = 259

$wb.Save()
